$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8200883333333334
$ws.Range("H2").Value = 2.460265
$ws.Range("I2").Value = 0.2405117342909232
$ws.Range("J2").Value = 0.2405117342909232
$ws.Range("M2").Value = 2.479667666666667
$ws.Range("N2").Value = 7.439003
$ws.Range("O2").Value = 0.02877351812554147
$ws.Range("P2").Value = 0.02877351812554146
$ws.Range("Q2").Value = 2.033546523977222
$ws.Range("R2").Value = 18.301918715795
$ws.Range("S2").Value = 0.006920368746025291
$ws.Range("T2").Value = 0.00692036874602529
$ws.Range("G3").Value = 0.8200883333333334
$ws.Range("H3").Value = 2.460265
$ws.Range("I3").Value = 0.2405117342909232
$ws.Range("J3").Value = 0.2405117342909232
$ws.Range("O3").Value = 0.8517172368296149
$ws.Range("P3").Value = 0.8517172368296148
$ws.Range("Q3").Value = 60.19446835835112
$ws.Range("R3").Value = 541.75021522516
$ws.Range("S3").Value = 0.2048479897553636
$ws.Range("T3").Value = 0.2048479897553636
$ws.Range("G4").Value = 0.8200883333333334
$ws.Range("H4").Value = 2.460265
$ws.Range("I4").Value = 0.2405117342909232
$ws.Range("J4").Value = 0.2405117342909232
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.497944666666667
$ws.Range("N4").Value = 10.493834
$ws.Range("O4").Value = 0.0405893804324885
$ws.Range("P4").Value = 0.0405893804324885
$ws.Range("Q4").Value = 2.868623611778889
$ws.Range("R4").Value = 25.81761250601
$ws.Range("S4").Value = 0.00976222228161187
$ws.Range("T4").Value = 0.00976222228161187
$ws.Range("G5").Value = 0.8200883333333334
$ws.Range("H5").Value = 2.460265
$ws.Range("I5").Value = 0.2405117342909232
$ws.Range("J5").Value = 0.2405117342909232
$ws.Range("M5").Value = 5.285651
$ws.Range("N5").Value = 15.856953
$ws.Range("O5").Value = 0.06133353146400924
$ws.Range("P5").Value = 0.06133353146400923
$ws.Range("Q5").Value = 4.334700719171667
$ws.Range("R5").Value = 39.012306472545
$ws.Range("S5").Value = 0.01475143402259576
$ws.Range("T5").Value = 0.01475143402259576
$ws.Range("G6").Value = 0.8200883333333334
$ws.Range("H6").Value = 2.460265
$ws.Range("I6").Value = 0.2405117342909232
$ws.Range("J6").Value = 0.2405117342909232
$ws.Range("M6").Value = 1.515569333333333
$ws.Range("N6").Value = 4.546708
$ws.Range("O6").Value = 0.01758633314834587
$ws.Range("P6").Value = 0.01758633314834587
$ws.Range("Q6").Value = 1.242900728624445
$ws.Range("R6").Value = 11.18610655762
$ws.Range("S6").Value = 0.004229719485326616
$ws.Range("T6").Value = 0.004229719485326616
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.589676
$ws.Range("H7").Value = 7.769028
$ws.Range("I7").Value = 0.7594882657090768
$ws.Range("J7").Value = 0.7594882657090768
$ws.Range("M7").Value = 2.479667666666667
$ws.Range("N7").Value = 7.439003
$ws.Range("O7").Value = 0.02877351812554147
$ws.Range("P7").Value = 0.02877351812554146
$ws.Range("Q7").Value = 6.421535844342666
$ws.Range("R7").Value = 57.793822599084
$ws.Range("S7").Value = 0.02185314937951617
$ws.Range("T7").Value = 0.02185314937951617
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.589676
$ws.Range("H8").Value = 7.769028
$ws.Range("I8").Value = 0.7594882657090768
$ws.Range("J8").Value = 0.7594882657090768
$ws.Range("O8").Value = 0.8517172368296149
$ws.Range("P8").Value = 0.8517172368296148
$ws.Range("Q8").Value = 190.0821700593813
$ws.Range("R8").Value = 1710.739530534432
$ws.Range("S8").Value = 0.6468692470742512
$ws.Range("T8").Value = 0.6468692470742511
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.589676
$ws.Range("H9").Value = 7.769028
$ws.Range("I9").Value = 0.7594882657090768
$ws.Range("J9").Value = 0.7594882657090768
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.497944666666667
$ws.Range("N9").Value = 10.493834
$ws.Range("O9").Value = 0.0405893804324885
$ws.Range("P9").Value = 0.0405893804324885
$ws.Range("Q9").Value = 9.058543352594667
$ws.Range("R9").Value = 81.52689017335199
$ws.Range("S9").Value = 0.03082715815087663
$ws.Range("T9").Value = 0.03082715815087663
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.589676
$ws.Range("H10").Value = 7.769028
$ws.Range("I10").Value = 0.7594882657090768
$ws.Range("J10").Value = 0.7594882657090768
$ws.Range("M10").Value = 5.285651
$ws.Range("N10").Value = 15.856953
$ws.Range("O10").Value = 0.06133353146400924
$ws.Range("P10").Value = 0.06133353146400923
$ws.Range("Q10").Value = 13.688123539076
$ws.Range("R10").Value = 123.193111851684
$ws.Range("S10").Value = 0.04658209744141347
$ws.Range("T10").Value = 0.04658209744141346
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.589676
$ws.Range("H11").Value = 7.769028
$ws.Range("I11").Value = 0.7594882657090768
$ws.Range("J11").Value = 0.7594882657090768
$ws.Range("M11").Value = 1.515569333333333
$ws.Range("N11").Value = 4.546708
$ws.Range("O11").Value = 0.01758633314834587
$ws.Range("P11").Value = 0.01758633314834587
$ws.Range("Q11").Value = 3.924833528869333
$ws.Range("R11").Value = 35.323501759824
$ws.Range("S11").Value = 0.01335661366301926
$ws.Range("T11").Value = 0.01335661366301925
